# Swap the weekly price data between row 2 and row 3.
# Columns affected: D (Fecha), J (Volumen), K (Precio mínimo),
# L (Precio máximo), M (Precio promedio ponderado), P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
